$wb = $excel.ActiveWorkbook

# Add a new worksheet right after "ODI Bowling" (i.e. at the end of the tab strip)
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item("ODI Bowling"))
$newSheet.Name = "ODI Batting Extra"

# --- Header row (row 1) ---
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Match the bold / bordered / centered header style used on the other sheets
$newSheet.Range("A1:F1").Font.Bold = $true
$newSheet.Range("A1:F1").HorizontalAlignment = -4108
$newSheet.Range("A1:F1").VerticalAlignment = -4160
$newSheet.Range("A1:F1").Borders.LineStyle = 1

# --- Data rows: force text-typed columns to stay as literal text (not get
# auto-coerced into numbers/percentages) by pre-setting the cell format ---
$textCells = "A2","C2","D2","E2","F2","A3","C3","D3","E3","F3"
foreach ($addr in $textCells) {
    $newSheet.Range($addr).NumberFormat = "@"
}

# Data row 2
$newSheet.Range("A2").Value = "4458"
$newSheet.Range("B2").Value = 5
$newSheet.Range("C2").Value = "0"
$newSheet.Range("D2").Value = "0"
$newSheet.Range("E2").Value = "1.09%"
$newSheet.Range("F2").Value = "NO"

# Data row 3
$newSheet.Range("A3").Value = "4459"
$newSheet.Range("B3").Value = 5
$newSheet.Range("C3").Value = "2"
$newSheet.Range("D3").Value = "0"
$newSheet.Range("E3").Value = "2.78%"
$newSheet.Range("F3").Value = "NO"

# Restore the originally active sheet/tab ("Player Info" was the active tab
# before this edit) so adding the sheet doesn't change workbook selection.
$wb.Worksheets.Item("Player Info").Activate()
